$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of vote data (row 31) - 10k+ votes timestamp for Balen
$ws.Range("A31").Value = "5/17/2022 8:30"
$ws.Range("A31").NumberFormat = "@"
$ws.Range("B31").Value = 10359
$ws.Range("C31").Value = 5515
$ws.Range("D31").Value = 4792
$ws.Range("E31").Formula = "=B31-C31"
$ws.Range("F31").Formula = "=B31-D31"

# Move the active selection like the source edit did
$ws.Range("F35").Select()
